# Fix bug: remove time T Z from brief vo
# Slide 4, "Content Placeholder 2" shape: split the last bullet's run so the
# trailing description reads "...服务器端" + "数据维护和界面" instead of
# "...服务器端数据更新和维护".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

$oldFragment = "数据更新和维护"
$newFragment = "数据维护和界面"

$fullText = $tr.Text
$idx = $fullText.IndexOf($oldFragment)

if ($idx -ge 0) {
    # TextRange/Characters is 1-indexed.
    $sub = $tr.Characters($idx + 1, $oldFragment.Length)
    $sub.Text = $newFragment
}
